$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark from its original location
#    (right before "Japan Product Manager Conference 2016" in the
#    contact-info paragraph near the end of the document).
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Insert "配信・" just before "撮影への同意" in the recording/
#    filming consent table cell, with its own distinct run
#    formatting (Arial Unicode MS / MS Mincho, 8pt).
# ------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = "撮影への同意"
$find1.Execute($null, $null, $null, $null, $null, $null, $true, 1, $null, "配信・撮影への同意", 2) | Out-Null

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Replacement.Font.Name = "Arial Unicode MS"
$find2.Replacement.Font.NameFarEast = "ＭＳ 明朝"
$find2.Replacement.Font.NameBi = "ＭＳ 明朝"
$find2.Text = "配信・"
$find2.Replacement.Text = "配信・"
$find2.Execute($null, $null, $null, $null, $null, $null, $true, 1, $null, $null, 2) | Out-Null

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark (collapsed) right at the new
#    insertion point, between "配信・" and "撮影への同意".
# ------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("撮影への同意", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPoint = $r3.Start
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4) Mark the "Slideshare" mention as a spell-checker exception span
#    (wraps it between proofErr spellStart/spellEnd markers).
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Slideshare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Write-Output "done"
